$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.412.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.892.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.88%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.15"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.691"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.20"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +9.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.351"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.01"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0737"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.27%  "
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.169.82"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.725"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.923.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.381.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "246.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.82"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -10.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.55"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.49"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("E30").Value = "  -3.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.128.40"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.80"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +9.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.852"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.55%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -19.61%  "
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.61"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0671"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("E43").Value = "  -2.10%  "
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.289.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.36"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0805"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.76%  "
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("E51").Value = "  -5.08%  "
